$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (used both in workbook.xml sheet name and tab)
$ws.Name = "Through 2022-03-31"

# Update the label for the March row (shared string)
$ws.Range("A4").Value = "March (through 03-31)"

# Update the March row values (row 4)
$ws.Range("C4").Value = 41
$ws.Range("D4").Value = 58
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 81
$ws.Range("I4").Value = 133

# Update the Total row values (row 5)
$ws.Range("C5").Value = 128
$ws.Range("D5").Value = 189
$ws.Range("F5").Value = 110
$ws.Range("G5").Value = 198
$ws.Range("H5").Value = 423
$ws.Range("I5").Value = 433
